$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values like
# "1.00" or "0.999" are not auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.748.40'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.035.90'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.18'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.69'
$ws.Range('E6').Value = '  +3.14%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.55'
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.111'
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  +4.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.544.62'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.78'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('E15').Value = '  +8.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.793.17'
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.28'
$ws.Range('E17').Value = '  +8.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.034.48'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.89'
$ws.Range('E19').Value = '  +3.35%  '
$ws.Range('E20').Value = '  +3.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '333.96'
$ws.Range('E21').Value = '  +2.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  +5.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.82'
$ws.Range('E24').Value = '  +4.05%  '
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0938'
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.84'
$ws.Range('E28').Value = '  +5.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.53'
$ws.Range('E29').Value = '  +8.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.82'
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.84'
$ws.Range('E32').Value = '  +1.37%  '
$ws.Range('E33').Value = '  +6.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.01'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.90'
$ws.Range('E35').Value = '  +5.68%  '
$ws.Range('E36').Value = '  +1.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0686'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.071.16'
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.52'
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.89'
$ws.Range('E41').Value = '  +8.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.319.87'
$ws.Range('E43').Value = '  +2.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.657'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  +5.04%  '
$ws.Range('E48').Value = '  +2.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.69'
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('E50').Value = '  -4.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0897'
